$wb = $excel.ActiveWorkbook

# Find the last existing sheet (2026-01-14) to insert the new sheet after it,
# and to use as a formatting template (header style + highlight style).
$sheetCount = $wb.Worksheets.Count
$templateSheet = $wb.Worksheets.Item($sheetCount)

# Add the new weekly sheet right after the last existing one.
$newSheet = $wb.Worksheets.Add($null, $templateSheet)
$newSheet.Name = "2026-01-21"

# Copy the header row formatting (bold, border, centered) from the template sheet.
$templateSheet.Range("A1:D1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "rank"
$newSheet.Range("B1").Value = "title"
$newSheet.Range("C1").Value = "volume"
$newSheet.Range("D1").Value = "publisher"

# Grab a cell that carries the "new/highlighted volume" fill from the template
# sheet so we can stamp the same style (rather than inventing a new one).
$highlightSource = $templateSheet.Range("C4")

$data = @(
  @(1, '転生貴族、鑑定スキルで成り上がる ~弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた~', 20, $false),
  @(2, '僕の心のヤバイやつ', 13, $false),
  @(3, 'ダンダダン', 22, $false),
  @(4, '弱虫ペダル', 98, $false),
  @(5, '十字架のろくにん', 23, $false),
  @(6, '貴族転生 ~恵まれた生まれから最強の力を得る~', 10, $false),
  @(7, '空挺ドラゴンズ', 21, $false),
  @(8, 'うるわしの宵の月', 10, $false),
  @(9, '呪術廻戦≡(モジュロ)', 1, $true),
  @(10, 'カグラバチ', 10, $false),
  @(11, '追放された没落令嬢は拳ひとつで異世界を生き延びる! コミック版', 1, $true),
  @(12, 'ヘルモード ~やり込み好きのゲーマーは廃設定の異世界で無双する~はじまりの召喚士13', 13, $false),
  @(13, '異世界転生で賢者になって冒険者生活 ~で異世界最強~', 11, $false),
  @(14, '領民0人スタートの辺境領主様 ~青のディアスと蒼角の乙女~14', 14, $false),
  @(15, '葬送のフリーレン', 15, $false),
  @(16, '七つ屋志のぶの宝石匣', 26, $false),
  @(17, 'あなたの寵妃でかまわない ~騎士令嬢は吸血公爵に溺愛される~ コミック版', 1, $true),
  @(18, '育成スキルはもういらないと勇者パーティを解雇されたので、退職金がわりにもらったを強くしてみる', 13, $false),
  @(19, '即死チートが最強すぎて、異世界のやつらがまるで相手にならないんですが。 -ΑΩ-13', 13, $false),
  @(20, '追放された転生重騎士はゲーム知識で無双する', 16, $false),
  @(21, '片田舎のおっさん、剣聖になる~ただの田舎の剣術師範だったのに、大成した弟子たちが俺を放ってくれない件~', 8, $false),
  @(22, 'オオカミ陛下の躾け方 悪役令嬢は破滅フラグと一夜を共にしてしまったけど、溺愛させます! コミック版', 1, $true),
  @(23, 'アラフォー賢者の異世界生活日記~気ままな異世界教師ライフ~', 17, $false),
  @(24, '脇役に転生したはずが、いつの間にか伝説の錬金術師になってた ~仲間たちが英雄でも俺は支援職なんだが~', 6, $false),
  @(25, '忘却バッテリー', 23, $false),
  @(26, '転生したらドラゴンの卵だった ~イバラのドラゴンロード~9', 9, $false),
  @(27, '恋せよまやかし天使ども', 6, $false),
  @(28, 'るろうに剣心―明治剣客浪漫譚・北海道編―', 10, $false),
  @(29, '攻略対象がモブ執事になりました', 1, $true),
  @(30, '馬小屋暮らしのご令嬢は案外領主に向いている? コミック版', 1, $true),
  @(31, 'ゲーム中盤で死ぬ悪役貴族に転生したので、外れスキルを駆使して最強を目指してみた', 4, $false),
  @(32, 'ダンジョン・シェルパ 迷宮道先案内人', 15, $false),
  @(33, '宮廷をクビになった植物魔導師はスローライフを謳歌する~のんびり世界樹を育てたら、最強領地ができました~', 5, $false),
  @(34, '映像研には手を出すな!', 10, $false),
  @(35, '侯爵令嬢は手駒を演じる 1(アリアンローズコミックス)', 1, $true),
  @(36, 'ハサミ男とサブカル女', 3, $true),
  @(37, '地雷なんですか?地原さん', 8, $false),
  @(38, '最強で最速の無限レベルアップ ~スキルとでレベル上限の枷が外れた俺は無双する~', 10, $false),
  @(39, 'ザ・ファブル The third secret', 3, $true),
  @(40, '漫画 ゆうえんち -バキ外伝-', 9, $false),
  @(41, 'ゴブリンスレイヤー', 17, $false),
  @(42, '最強タンクの迷宮攻略 ~体力9999のレアスキル持ちタンク、勇者パーティーを追放される~', 14, $false),
  @(43, '冒険家になろう!~スキルボードでダンジョン攻略~(コミック)', 12, $false),
  @(44, '二十と成獣', 7, $false),
  @(45, '僕の心のヤバイやつ ラブコメディが始まらない', 1, $true),
  @(46, 'あかね噺', 20, $false),
  @(47, 'くちべた食堂', 1, $true),
  @(48, '堕天使そぷらのちゃんの復讐', 1, $true),
  @(49, '咲きそめコンプレックス', 1, $true),
  @(50, 'ナイト・リセット・キロポスト', 1, $true),
  @(51, 'エモロイド', 1, $true),
  @(52, '朝比奈くんは一途に愛したい', 9, $false),
  @(53, 'このたび鬼上司の秘書になりまして', 6, $false),
  @(54, '黒猫男子は年上彼女を溺愛する', 3, $true),
  @(55, 'パラレルトラッパーズ!', 1, $true),
  @(56, '起きたら20年後なんですけど! ~悪役令嬢のその後のその後~ 1(アリアンローズコミックス)', 1, $true),
  @(57, 'お腹が痛くて死にそうです ~夜職女子の潰瘍性大腸炎日記~1', 1, $true),
  @(58, '不倫がバレて謝罪代行を使いました。 不倫がバレて謝罪代行を使いました。', 1, $true),
  @(59, '雷雷雷', 6, $false),
  @(60, 'ザ・ファブル The third secret', 1, $true),
  @(61, '薬屋のひとりごと', 16, $false),
  @(62, '穏やか貴族の休暇のすすめ。@COMIC', 14, $false),
  @(63, '冒険王ビィト', 19, $false),
  @(64, '呪術廻戦', 29, $false),
  @(65, 'くちべた食堂', 2, $true),
  @(66, 'くちべた食堂', 3, $true),
  @(67, 'スレイマンズ', 1, $true),
  @(68, '堕天使そぷらのちゃんの復讐', 2, $true),
  @(69, '堕天使そぷらのちゃんの復讐', 3, $true),
  @(70, 'ディメンションウェーブ', 1, $true),
  @(71, '咲きそめコンプレックス', 2, $true),
  @(72, '咲きそめコンプレックス', 3, $true),
  @(73, 'ナイト・リセット・キロポスト', 2, $true),
  @(74, 'ナイト・リセット・キロポスト', 3, $true),
  @(75, '幸福のおいしい道すがら', 1, $true),
  @(76, '春川さんは今日も飢えている', 1, $true),
  @(77, '年の差十五の旦那様~辺境伯の花嫁候補~', 1, $true),
  @(78, '年の差十五の旦那様~辺境伯の花嫁候補~', 2, $true),
  @(79, '年の差十五の旦那様~辺境伯の花嫁候補~', 3, $true),
  @(80, '年の差十五の旦那様~辺境伯の花嫁候補~', 4, $false),
  @(81, '年の差十五の旦那様~辺境伯の花嫁候補~', 5, $false),
  @(82, '年の差十五の旦那様~辺境伯の花嫁候補~', 6, $false),
  @(83, '前世は保育士、今世は悪役令嬢?からの、わがまま姫様の教育係!?', 1, $true),
  @(84, '高嶺の花宮くんとぼっちな彼女', 1, $true),
  @(85, '勇者の出番ねぇからっ!! ~異世界転生するけど俺は脇役と言われました~ コミック版', 1, $true),
  @(86, '朝比奈くんは一途に愛したい', 6, $false),
  @(87, '朝比奈くんは一途に愛したい', 7, $false),
  @(88, '朝比奈くんは一途に愛したい', 8, $false),
  @(89, 'このたび鬼上司の秘書になりまして', 5, $false),
  @(90, '元カレと再会してハジメテ。10年分抱かせて', 3, $true),
  @(91, '元カレと再会してハジメテ。10年分抱かせて', 4, $false),
  @(92, '社長が私を抱く理由', 3, $true),
  @(93, '社長が私を抱く理由', 4, $false),
  @(94, 'ひみつの犬飼くん', 3, $true),
  @(95, '乙女ゲー世界はモブに厳しい世界です', 1, $true),
  @(96, '外来魔法生物対策課', 1, $true),
  @(97, 'パラレルトラッパーズ!', 2, $true),
  @(98, 'パラレルトラッパーズ!', 3, $true),
  @(99, '男爵無双 貴族嫌いの青年が田舎貴族に転生した件', 1, $true),
  @(100, '魔法学園の最強暗殺者', 1, $true)
)

$row = 2
foreach ($item in $data) {
    $rank = $item[0]
    $title = $item[1]
    $volume = $item[2]
    $isHighlighted = $item[3]

    $newSheet.Cells.Item($row, 1).Value = $rank
    $newSheet.Cells.Item($row, 2).Value = $title
    $newSheet.Cells.Item($row, 3).Value = $volume

    if ($isHighlighted) {
        $highlightSource.Copy()
        $newSheet.Cells.Item($row, 3).PasteSpecial(-4122)
        $newSheet.Cells.Item($row, 3).Value = $volume
    }

    $row = $row + 1
}

Write-Output "Added sheet 2026-01-21 with $($data.Count) ranking rows"
